$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 52, pushing existing rows 52:91 down to 53:92
$ws.Rows(52).Insert()

# Populate the new row 52 with its data. Columns that are constant across
# the whole "Achicoria" block (A,B,C,E,F,G,H,I,N,O,Q,R) are copied from the
# row immediately below (old row 52, now shifted to row 53).
$ws.Range("A52").Value = 9
$ws.Range("B52").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C52").Value = "Metropolitana"
$ws.Range("D52").Value = 45233
$ws.Range("E52").Value = 13
$ws.Range("F52").Value = 100112010
$ws.Range("G52").Value = "Achicoria"
$ws.Range("H52").Value = "Sin especificar"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 70
$ws.Range("K52").Value = 8000
$ws.Range("L52").Value = 8000
$ws.Range("M52").Value = 8000
$ws.Range("N52").Value = '$/caja 16 unidades'
$ws.Range("O52").Value = "Provincia de Quillota"
$ws.Range("P52").Value = 500
$ws.Range("Q52").Value = 16
$ws.Range("R52").Value = "Hortaliza"
